$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (and two row coin-swaps)
# Cells whose new text would otherwise be auto-parsed as a number by Excel
# are temporarily forced to Text format, then the format is cleared again
# afterwards so no extra style survives in the saved workbook.

$ws.Range("D2").Value = "28.224.71"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "1.795.68"
$ws.Range("E3").Value = "  +2.08%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.37"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4521"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +15.80%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3743"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +9.95%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.76"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.149"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07550"
$ws.Range("D11").ClearFormats()
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.57"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.000"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.307"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.569"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +6.70%  "
$ws.Range("D16").Value = "1.792.67"
$ws.Range("E16").Value = "  +1.99%  "
$ws.Range("E17").Value = "  +3.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06729"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "81.11"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.87%  "
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.55"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.371"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.46%  "
$ws.Range("D23").Value = "28.243.89"
$ws.Range("E23").Value = "  +0.79%  "
$ws.Range("E24").Value = "  +1.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.421"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.60"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +3.13%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.361"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.26%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.57"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.99%  "
$ws.Range("D29").Value = "1.996.49"
$ws.Range("E29").Value = "  +1.94%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.24"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.236"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.021"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.23%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.827"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.09440"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +8.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.2347"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +11.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.17"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06343"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.90%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02335"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.198"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6585"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.364"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +6.14%  "
$ws.Range("E42").Value = "  -1.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.212"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.60%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.12"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.91%  "
$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9996"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6108"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.86%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.796"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.94"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.030"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.85%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07129"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.163"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.61%  "
